# update doc presentations per hyperlinks
#
# 1) Slide 1 title: merge the separate "Drug " / "Prescriptions " / "& " /
#    "Reactions " runs into a single run "Drug Prescriptions & Reactions ".
# 2) Slide 8 body: merge "4970 " + "Visible Nodes (Drug Matrix)" and
#    "2908 " + "Hidden Nodes (Reactions)" into single runs per paragraph.
# 3) Slide 9 title: merge "Top " + "Similar Drugs Sample" into a single run.

$p = $ppt.ActivePresentation

# --- Slide 1: title "OpenFDA Adverse<br>Drug Prescriptions & Reactions <br>Prediction Model"
$s1 = $p.Slides.Item(1)
$titleShape = $s1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$fullTitle = $titleRange.Text
$start = $fullTitle.IndexOf("Drug Prescriptions & Reactions ") + 1
$len = "Drug Prescriptions & Reactions ".Length
$titleRange.Characters($start, $len).Text = "Drug Prescriptions & Reactions "

# --- Slide 8: stats body "4970 Visible Nodes (Drug Matrix)" / "2908 Hidden Nodes (Reactions)"
$s8 = $p.Slides.Item(8)
$statsShape = $s8.Shapes.Item(2)
$statsRange = $statsShape.TextFrame.TextRange

$statsText = $statsRange.Text
$start = $statsText.IndexOf("4970 Visible Nodes (Drug Matrix)") + 1
$len = "4970 Visible Nodes (Drug Matrix)".Length
$statsRange.Characters($start, $len).Text = "4970 Visible Nodes (Drug Matrix)"

$statsText = $statsRange.Text
$start = $statsText.IndexOf("2908 Hidden Nodes (Reactions)") + 1
$len = "2908 Hidden Nodes (Reactions)".Length
$statsRange.Characters($start, $len).Text = "2908 Hidden Nodes (Reactions)"

# --- Slide 9: title "Top Similar Drugs Sample"
$s9 = $p.Slides.Item(9)
$sampleTitleShape = $s9.Shapes.Item(1)
$sampleTitleRange = $sampleTitleShape.TextFrame.TextRange
$sampleText = $sampleTitleRange.Text
$start = $sampleText.IndexOf("Top Similar Drugs Sample") + 1
$len = "Top Similar Drugs Sample".Length
$sampleTitleRange.Characters($start, $len).Text = "Top Similar Drugs Sample"
